$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("w01_100.5_optimization_results")

$data = @(
    @(1, 608, 35.57893896159422, 2, -1, 1.001),
    @(2, 613, 14.43370405830502, 0.3079846385005394, -0.8758331993145428, 2.642459626525638),
    @(3, 617, 2.472059333903036, 1.338187291993382, -0.5689140721391157, 1.981702082060714),
    @(4, 633, 37.46487104357716, 2, -1, 3),
    @(5, 619, 35.60018327439182, 1.727771941716866, -0.06507520523356825, 1.200327548986884),
    @(6, 621, 33.55139664766729, 1.194904044790017, -0.6280816901515246, 1.729623129306444),
    @(7, 633, 10.94186633782868, 1.204772014142755, -0.9282523303781325, 2.271917763974799),
    @(8, 627, 5.269805997601083, 0.5943730024893901, -0.8790064180212761, 1.415949102927504),
    @(9, 643, 34.54245327933843, 1.876331964609189, -0.9408624466819731, 2.88886469693255),
    @(10, 623, 16.67757733126458, 0.2469772544864867, -1, 1.447108540977185)
)

$startRow = 12
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
